$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-numeric-looking Price (D) values to stay as Text, matching the
# original inlineStr storage, by pre-setting NumberFormat to "@" before assignment.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D15", "D20", "D22", "D23", "D24", "D36", "D40", "D42", "D47", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.932.88"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").Value = "1.638.14"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").Value = "214.65"

$ws.Range("D6").Value = "0.506"
$ws.Range("E6").Value = "  +0.62%  "

$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "0.0637"
$ws.Range("E8").Value = "  +0.69%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.254"
$ws.Range("E9").Value = "  -1.04%  "

$ws.Range("D10").Value = "19.62"
$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").Value = "1.863.89"

$ws.Range("E13").Value = "  -0.50%  "

$ws.Range("D14").Value = "1.655.44"
$ws.Range("E14").Value = "  +1.07%  "

$ws.Range("D15").Value = "0.544"
$ws.Range("E15").Value = "  -1.63%  "

$ws.Range("D16").Value = "0.0₃0757"
$ws.Range("E16").Value = "  -0.66%  "

$ws.Range("E17").Value = "  -0.87%  "

$ws.Range("D18").Value = "25.951.15"
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("D20").Value = "193.86"
$ws.Range("E20").Value = "  +0.74%  "

$ws.Range("E21").Value = "  -1.34%  "

$ws.Range("D22").Value = "9.90"
$ws.Range("E22").Value = "  -0.85%  "

$ws.Range("D23").Value = "6.27"
$ws.Range("E23").Value = "  -1.25%  "

$ws.Range("D24").Value = "1.81"
$ws.Range("E24").Value = "  +0.43%  "

$ws.Range("E25").Value = "  +0.88%  "

$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("E27").Value = "  +2.92%  "

$ws.Range("E28").Value = "  -0.53%  "

$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("E31").Value = "  +1.50%  "

$ws.Range("E32").Value = "  -1.12%  "

$ws.Range("E33").Value = "  -1.04%  "

$ws.Range("E34").Value = "  -2.77%  "

$ws.Range("E35").Value = "  +1.30%  "

$ws.Range("D36").Value = "0.903"
$ws.Range("E36").Value = "  -0.79%  "

$ws.Range("D37").Value = "1.138.09"
$ws.Range("E37").Value = "  -0.97%  "

$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("E39").Value = "  -1.80%  "

$ws.Range("D40").Value = "0.0157"
$ws.Range("E40").Value = "  +0.08%  "

$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("D42").Value = "99.42"
$ws.Range("E42").Value = "  -1.30%  "

$ws.Range("E43").Value = "  -0.99%  "

$ws.Range("E44").Value = "  -4.14%  "

$ws.Range("D45").Value = "1.773.03"
$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("D46").Value = "0.0₆0114"
$ws.Range("E46").Value = "  +2.09%  "

$ws.Range("D47").Value = "56.49"
$ws.Range("E47").Value = "  +1.41%  "

$ws.Range("E48").Value = "  +3.66%  "

$ws.Range("E49").Value = "  -0.88%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.415"
$ws.Range("E50").Value = "  -0.49%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.63"
$ws.Range("E51").Value = "  +0.31%  "
